# Apply updated Betfair odds values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 (Manisa FK vs Umraniyespor)
$ws.Range("G5").Value = 1.86
$ws.Range("K5").Value = 4.7
$ws.Range("P5").Value = 2.28
$ws.Range("Q5").Value = 1.63

# Row 6 (Genoa vs Cagliari)
$ws.Range("AM6").Value = 190
$ws.Range("AN6").Value = 27

# Row 7 (Al-Ettifaq vs Al-Khaleej Saihat)
$ws.Range("P7").Value = 2.04

# Row 8 (Al-Hilal vs Al Nassr)
$ws.Range("G8").Value = 2.78
$ws.Range("H8").Value = 2.66
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 3.65
$ws.Range("K8").Value = 4.9
$ws.Range("P8").Value = 3.1
$ws.Range("Q8").Value = 1.33

# Row 10 (Huesca vs Cordoba)
$ws.Range("P10").Value = 1.62

# Row 11 (Juventus vs US Cremonese)
$ws.Range("F11").Value = 1.27
$ws.Range("G11").Value = 1.28
$ws.Range("X11").Value = 22
$ws.Range("AF11").Value = 7.2

# Row 12 (Sevilla vs Celta Vigo)
$ws.Range("G12").Value = 2.86
$ws.Range("H12").Value = 2.78

$wb.Save()
